$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: correct the evening session time range and recompute hours ---
$ws.Range("B9").Value = "18.00-20.30"
$ws.Range("G9").Value = 2.5

# --- Row 11: log the hours spent that day (previously left blank) ---
$ws.Range("G11").Value = 1

# --- Row 17: new diary entry, 28 loka (D'Alembert / Hooke's law, ideal stiff spring forcegenerator) ---
$ws.Range("A17").Value = "28 loka"

$ws.Range("B17").Value = "8.45-10.45, 11.45-13.15"
$ws.Range("B17").NumberFormat = $ws.Range("B16").NumberFormat
$ws.Range("B17").WrapText = $true

$ws.Range("C17").Value = "D'Alamembertin laki, Hookin laki, Oppikirjasta 79-103"
$ws.Range("C17").WrapText = $true

$ws.Range("D17").Value = "Noniin, nyt roppakaupalla uutta asiaa. Kirja on aika hyvä, siellä pedataan paljon että miksi asiat tehdään sillä tavalla ja mitä jatkossa luvassa."
$ws.Range("D17").WrapText = $true

$ws.Range("E17").Value = "Katsotaan mitä modernisointia tähän demoon löytyy, nyt kopioitu lähinnä boilerplatea. Ehkä huomenna koodaillaan enemmän."
$ws.Range("E17").WrapText = $true

$ws.Range("G17").Value = 3.5

$ws.Rows.Item(17).RowHeight = 58

# --- Update the running total to cover the extended range ---
$ws.Range("H3").Formula = "=SUM(G3:G60)"
$ws.Range("H3").NumberFormat = $ws.Range("G3").NumberFormat

# --- Column H a touch wider to fit the new total ---
$ws.Range("H1").ColumnWidth = 9.43

# --- Selection as left by the author ---
$ws.Range("H12").Select()
